# Update ObjTables/SBtab header strings: bump version 0.0.9 -> 1.0.0 and date to 2020-05-29 00:21:23
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()
$ws.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='1.0.0' date='2020-05-29 00:21:23'"
$ws.Range("A2").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Compartment' name='Compartment' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(2)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Compound' name='Compound' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(3)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Definition' name='Definition' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(4)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Enzyme' name='Enzyme' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(5)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='FbcObjective' name='FbcObjective' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(6)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Gene' name='Gene' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(7)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Layout' name='Layout' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(8)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Measurement' name='Measurement' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(9)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='PbConfig' name='PbConfig' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(10)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Position' name='Position' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(11)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Protein' name='Protein' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(12)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Quantity' name='Quantity' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(13)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='QuantityInfo' name='QuantityInfo' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(14)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='QuantityMatrix' name='QuantityMatrix' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(15)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Reaction' name='Reaction' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(16)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(17)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Regulator' name='Regulator' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(18)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Relation' name='Relation' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(19)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Relationship' name='Relationship' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(20)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrix' name='SparseMatrix' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(21)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(22)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(23)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrixRow' name='SparseMatrixRow' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(24)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(25)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='rxnconContingencyList' name='rxnconContingencyList' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item(26)
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='rxnconReactionList' name='rxnconReactionList' date='2020-05-29 00:21:23' objTablesVersion='1.0.0'"
$ws.Protect()
